# Online Enrollment Check List - Aug 21 executive-meeting pass for Colorado
#
# Insert 7 new rows right above the existing "NANNY" section (old row 51),
# pushing everything from there down by 7 rows, and use the 7 freshly
# inserted rows to record 3 new to-do notes (plus 4 blank spacer rows)
# directly under the "Colorado contract" block that ends at row 50.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 7 blank rows before the old row 51 ("NANNY" header), shifting the
# rest of the sheet down. Excel carries the formatting of the row above
# (row 50, style s="4") onto the newly inserted rows, matching the target.
$ws.Rows("51:57").Insert()

# Fill in the three new notes in column D of the freshly inserted rows.
$ws.Range("D51").Value = "Need to add Ifee to items being sent to database, need to get the upc code"
$ws.Range("D52").Value = "Kids names need to be upper case"
$ws.Range("D53").Value = "Remove just Dues from colorado contract and show total only"

# Rows D54:D57 stay blank (already created blank by the row insert above).

# Match the author's saved view state: scrolled down with D53 selected.
$ws.Range("D53").Select()
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1
